$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7038
$ws1.Range("F7").Value = 153
$ws1.Range("F8").Value = 119
$ws1.Range("F10").Value = 13
$ws1.Range("F17").Value = 3634
$ws1.Range("G17").Value = 58.5
$ws1.Range("F23").Value = 2264
$ws1.Range("F32").Value = 272
$ws1.Range("F33").Value = 105

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7038
$ws4.Range("F8").Value = 153
$ws4.Range("F9").Value = 119
$ws4.Range("F11").Value = 13
$ws4.Range("F18").Value = 3634
$ws4.Range("G18").Value = 58.5
$ws4.Range("F24").Value = 2264
$ws4.Range("F33").Value = 272
$ws4.Range("F34").Value = 105
